# The page-XML for this transcription had each "<id>...</id>" marker
# split into several runs (<id> / the id text, sometimes itself split
# into a few runs / </id>) with inconsistent leftover formatting from
# the page it was copy/pasted from. Tidy that up: each "<id>...</id>"
# becomes a single run.
#
# This touches every paragraph in the doc whose text is exactly
# "<id>XXXX</id>" - currently p050r_1 and p050r_2.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $pRange = $para.Range
    $text = $pRange.Text

    # Paragraph.Range.Text includes the trailing paragraph mark; strip it
    # (and any other trailing control chars) before comparing.
    $trimmed = $text.TrimEnd([char]13, [char]7)

    if ($trimmed.StartsWith("<id>") -and $trimmed.EndsWith("</id>")) {
        # Re-use the whole paragraph Range (start .. end, i.e. including
        # the trailing paragraph-mark slot). Assigning .Text against a
        # range like this collapses every run it spans into a single
        # run, which picks up the formatting of the range's first
        # character - i.e. the original "<id>" run's formatting
        # (Courier New / 7F6000 / 9pt), exactly what the merged run
        # should use.
        $idRange = $d.Range($pRange.Start, $pRange.End)
        $idRange.Text = $trimmed
    }
}
